$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.163.50'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '3.148.20'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.68'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.90'
$ws.Range("E6").Value = '  -2.80%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '3.138.68'
$ws.Range("E8").Value = '  +0.89%  '
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.26'
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.19'
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("D15").Value = '3.659.97'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '3.139.50'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '63.091.42'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.67'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '473.64'
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.10'
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.699'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.67'
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.89'
$ws.Range("E24").Value = '  -2.49%  '
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.02'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.97'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.12'
$ws.Range("E30").Value = '  +3.80%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.88'
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("E33").Value = '  -4.27%  '
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("E35").Value = '  -1.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.79'
$ws.Range("E36").Value = '  -2.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.28'
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("D38").Value = '0.0₃0692'
$ws.Range("E38").Value = '  -6.98%  '
$ws.Range("E39").Value = '  -0.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '419.95'
$ws.Range("E40").Value = '  -3.02%  '
$ws.Range("E41").Value = '  -6.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.21'
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("D43").Value = '2.912.83'
$ws.Range("E43").Value = '  +1.96%  '
$ws.Range("E44").Value = '  -5.43%  '
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  -6.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.37'
$ws.Range("E51").Value = '  -0.60%  '
